# Add industry ("DMD_IND-*") sets/tables to the VEDA_Sets-Proc sheet, and
# refresh the cell-selection / zoom view state that Excel re-stamps on save.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) VEDA_Sets-Comm: view state only (selection moved from D9 to D7, and
#    normal-view zoom pinned at 100%). No cell data changes on this sheet.
# ---------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item("VEDA_Sets-Comm")
[void]$wsComm.Activate()
$excel.ActiveWindow.Zoom = 100
[void]$wsComm.Range("D7").Select()

# ---------------------------------------------------------------------
# 2) VEDA_Sets-Proc: thirteen new "Industry" demand sets appended below
#    the existing table (rows 50-62), plus the updated selection.
# ---------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("VEDA_Sets-Proc")
[void]$wsProc.Activate()

# PSET_PN, SetName, SetDesc triples - one new row per industry sub-sector.
$industrySets = @(
    @("I*CAF*", "DMD_IND-CAF", "Industry - Chemicals and man-made fibers"),
    @("I*EOE*", "DMD_IND-EOE", "Industry - Electrical and optical equipment"),
    @("I*FAB*", "DMD_IND-FAB", "Industry - Food and beverages"),
    @("I*MAE*", "DMD_IND-MAE", "Industry - Machinery and equipment n.e.c."),
    @("I*MAP*", "DMD_IND-MAP", "Industry - Basic metals and fabricated metal prod."),
    @("I*NEM*", "DMD_IND-NEM", "Industry - Non-energy mining"),
    @("I*OMA*", "DMD_IND-OMA", "Industry - Other manufacturing"),
    @("I*ONM*", "DMD_IND-ONM", "Industry - Other non-metalic mineral products"),
    @("I*PX4*", "DMD_IND-PX4", "Industry - Pulp, paper, publishing and printing"),
    @("I*RAP*", "DMD_IND-RAP", "Industry - Rubber and plastic products"),
    @("I*TAP*", "DMD_IND-TAP", "Industry - Textiles and textile products"),
    @("I*TEM*", "DMD_IND-TEM", "Industry - Transport equipment manufacture"),
    @("I*WAP*", "DMD_IND-WAP", "Industry - Wood and wood products")
)

$startRow = 50
for ($i = 0; $i -lt $industrySets.Count; $i++) {
    $r = $startRow + $i
    $entry = $industrySets[$i]

    $wsProc.Cells.Item($r, 1).Value  = "DMD"      # A: PSET_SET
    $wsProc.Cells.Item($r, 2).Value  = $entry[0]  # B: PSET_PN
    $wsProc.Cells.Item($r, 6).Value  = $entry[1]  # F: SetName
    $wsProc.Cells.Item($r, 7).Value  = $entry[2]  # G: SetDesc
    $wsProc.Cells.Item($r, 8).Value  = "AND"      # H: T_Pos_AndOr
    $wsProc.Cells.Item($r, 9).Value  = "OR"       # I: T_Neg_AndOr
    $wsProc.Cells.Item($r, 10).Value = "AND"      # J: T_Pos_AndOr_forSets
    $wsProc.Cells.Item($r, 11).Value = "OR"       # K: T_Neg_AndOr_forSets
}

[void]$wsProc.Range("F65").Select()
